{"js": "// Replace the two-digit-number-divided-by-one-digit-number problems\n// in the table with the new set of problems, per the commit diff.\nconst replacements = [\n  [\"39\u00f79=\", \"64\u00f78=\"],\n  [\"26\u00f78=\", \"66\u00f72=\"],\n  [\"17\u00f73=\", \"52\u00f72=\"],\n  [\"49\u00f72=\", \"97\u00f76=\"],\n  [\"45\u00f78=\", \"24\u00f74=\"],\n  [\"36\u00f74=\", \"38\u00f73=\"],\n  [\"73\u00f78=\", \"21\u00f74=\"],\n  [\"40\u00f76=\", \"75\u00f77=\"],\n  [\"56\u00f77=\", \"50\u00f73=\"],\n  [\"42\u00f78=\", \"34\u00f73=\"],\n  [\"46\u00f74=\", \"74\u00f74=\"],\n  [\"17\u00f75=\", \"57\u00f79=\"],\n  [\"31\u00f78=\", \"58\u00f77=\"],\n  [\"51\u00f75=\", \"87\u00f77=\"],\n  [\"20\u00f78=\", \"61\u00f75=\"],\n  [\"43\u00f73=\", \"82\u00f79=\"],\n  [\"80\u00f79=\", \"78\u00f75=\"],\n  [\"68\u00f74=\", \"76\u00f73=\"],\n  [\"24\u00f73=\", \"49\u00f79=\"],\n  [\"87\u00f75=\", \"81\u00f79=\"],\n  [\"56\u00f78=\", \"67\u00f72=\"],\n  [\"49\u00f78=\", \"45\u00f76=\"],\n  [\"64\u00f74=\", \"93\u00f75=\"],\n  [\"76\u00f76=\", \"90\u00f79=\"],\n  [\"72\u00f79=\", \"24\u00f77=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the two-digit-number-divided-by-one-digit-number problems\n# in the table with the new set of problems, per the commit diff.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"39\u00f79=\", \"64\u00f78=\"),\n    @(\"26\u00f78=\", \"66\u00f72=\"),\n    @(\"17\u00f73=\", \"52\u00f72=\"),\n    @(\"49\u00f72=\", \"97\u00f76=\"),\n    @(\"45\u00f78=\", \"24\u00f74=\"),\n    @(\"36\u00f74=\", \"38\u00f73=\"),\n    @(\"73\u00f78=\", \"21\u00f74=\"),\n    @(\"40\u00f76=\", \"75\u00f77=\"),\n    @(\"56\u00f77=\", \"50\u00f73=\"),\n    @(\"42\u00f78=\", \"34\u00f73=\"),\n    @(\"46\u00f74=\", \"74\u00f74=\"),\n    @(\"17\u00f75=\", \"57\u00f79=\"),\n    @(\"31\u00f78=\", \"58\u00f77=\"),\n    @(\"51\u00f75=\", \"87\u00f77=\"),\n    @(\"20\u00f78=\", \"61\u00f75=\"),\n    @(\"43\u00f73=\", \"82\u00f79=\"),\n    @(\"80\u00f79=\", \"78\u00f75=\"),\n    @(\"68\u00f74=\", \"76\u00f73=\"),\n    @(\"24\u00f73=\", \"49\u00f79=\"),\n    @(\"87\u00f75=\", \"81\u00f79=\"),\n    @(\"56\u00f78=\", \"67\u00f72=\"),\n    @(\"49\u00f78=\", \"45\u00f76=\"),\n    @(\"64\u00f74=\", \"93\u00f75=\"),\n    @(\"76\u00f76=\", \"90\u00f79=\"),\n    @(\"72\u00f79=\", \"24\u00f77=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1  # wdFindContinue\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute($find.Text, $find.MatchCase, $find.MatchWholeWord, $find.MatchWildcards, $false, $false, $find.Forward, $find.Wrap, $false, $newText, 2) | Out-Null\n}\n\n$d.Save()\n"}
